# [ANV] updated decay chains
# Adds a new "Middle (frac): Hang's Normalization" column (I) to the
# "Target Fractions" sheet, computed as F{row} * 'Shotcrete Hang'!$B$13.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Target Fractions")

# New column width
$ws.Columns.Item(9).ColumnWidth = 31

# New column header (I1) - mirrors the format of the other header cells (H1)
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(1, 9).Value = "Middle (frac): Hang's Normalization"

# Rows that carry the new ratio formula in column I
$formulaRows = @(4, 6, 7, 9, 11, 12, 13, 15, 17, 18, 19, 21, 22, 23, 25, 26, 27, 28, 29, 31, 33, 34, 35, 36)
foreach ($r in $formulaRows) {
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 9).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 9).Formula = "=F$r*'Shotcrete Hang'!`$B`$13"
}

# Section-header rows: column I stays empty but shares the same numeric
# style as the other empty cells in F:H on those rows
$emptyRows = @(5, 8, 10, 14, 16, 20, 24, 30, 32)
foreach ($r in $emptyRows) {
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 9).PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

# Matching view-state tweak on the "Shotcrete Hang" sheet
$wsHang = $wb.Worksheets.Item("Shotcrete Hang")
$wsHang.Range("B2").Select()

# Update the stored selection to match the authored view state; this sheet
# ("Target Fractions") remains the active tab, so select it last
$ws.Activate()
$ws.Range("I32").Select()
